# Update column F ("dSF") values on the active worksheet to match the
# re-pulled / re-pushed data and refreshed mean calculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = -4
    5  = 8
    9  = -2
    10 = -7
    12 = -1
    14 = -1
    15 = -2
    20 = -1
    24 = 4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
